$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.104095339775085
$ws.Range("B1").Value = 1.707998871803284
$ws.Range("C1").Value = 6.882572650909424
$ws.Range("D1").Value = 1.629111647605896
$ws.Range("E1").Value = 0.9437104463577271
